# Apply the "cryptos list" refresh captured in the commit:
#   "Updated cryptos list on Wed Aug  7 23:37:34 UTC 2024 with GitHub Actions"
# Updates Price (D) / Volume(1h) (E) figures for each coin row, and fixes a
# mis-ordering of the VeChain / Bittensor rows (47 & 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.072.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.342.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.40%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "475.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +20.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.347.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0956"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.324"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.752.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.100.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.75%  "

$ws.Range("E17").Value = "  -4.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.347.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "313.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.88%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.393"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.152"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.450.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.10%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0739"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "145.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.09"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("E36").Value = "  -4.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.804"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.18%  "

$ws.Range("E39").Value = "  +10.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.576"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0514"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "248.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.46%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0220"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.790.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.70%  "

